$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated bus voltage magnitudes (vm_pu) for the 380 kV case, rows 2-25
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053896268566671
$ws.Range("D2").Value = 1.067881951492595
$ws.Range("E2").Value = 1.061465396656779
$ws.Range("F2").Value = 1.074476957839882
$ws.Range("I2").Value = 1.053606814279874
$ws.Range("J2").Value = 1.058911207378102
$ws.Range("K2").Value = 1.070588624021155
$ws.Range("L2").Value = 1.064189383248673
$ws.Range("M2").Value = 1.077166067730351
$ws.Range("N2").Value = 1.023343406525857
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.055001702244641
$ws.Range("D3").Value = 1.068708921242147
$ws.Range("E3").Value = 1.062431308224346
$ws.Range("F3").Value = 1.075482472832595
$ws.Range("I3").Value = 1.053950241372019
$ws.Range("J3").Value = 1.059666950345578
$ws.Range("K3").Value = 1.071231437012455
$ws.Range("L3").Value = 1.064969539885895
$ws.Range("M3").Value = 1.077988252558492
$ws.Range("N3").Value = 1.023602519392131
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.055717159302028
$ws.Range("D4").Value = 1.069244022243975
$ws.Range("E4").Value = 1.063056760753208
$ws.Range("F4").Value = 1.076133573261865
$ws.Range("I4").Value = 1.054171187091721
$ws.Range("J4").Value = 1.060155552370791
$ws.Range("K4").Value = 1.071646724590362
$ws.Range("L4").Value = 1.06547417378598
$ws.Range("M4").Value = 1.078520108698895
$ws.Range("N4").Value = 1.023769866794241
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056017977912843
$ws.Range("D5").Value = 1.069468977155693
$ws.Range("E5").Value = 1.06331980674236
$ws.Range("F5").Value = 1.076407406749893
$ws.Range("I5").Value = 1.054263767365754
$ws.Range("J5").Value = 1.060360861406719
$ws.Range("K5").Value = 1.071821154166299
$ws.Range("L5").Value = 1.065686278656967
$ws.Range("M5").Value = 1.078743664177515
$ws.Range("N5").Value = 1.023840143861906
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056068489064443
$ws.Range("D6").Value = 1.069506747965225
$ws.Range("E6").Value = 1.063363979512368
$ws.Range("F6").Value = 1.076453391085666
$ws.Range("I6").Value = 1.054279294092875
$ws.Range("J6").Value = 1.060395327874222
$ws.Range("K6").Value = 1.071850432422146
$ws.Range("L6").Value = 1.065721889457251
$ws.Range("M6").Value = 1.078781197955362
$ws.Range("N6").Value = 1.023851939242922
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.055721178695579
$ws.Range("D7").Value = 1.06924702810943
$ws.Range("E7").Value = 1.063060275171509
$ws.Range("F7").Value = 1.076137231803266
$ws.Range("I7").Value = 1.054172425353732
$ws.Range("J7").Value = 1.060158296110296
$ws.Range("K7").Value = 1.071649055945581
$ws.Range("L7").Value = 1.065477008111749
$ws.Range("M7").Value = 1.07852309600335
$ws.Range("N7").Value = 1.023770806137797
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.054269820572442
$ws.Range("D8").Value = 1.06816143016738
$ws.Range("E8").Value = 1.061791739026678
$ws.Range("F8").Value = 1.074816679787786
$ws.Range("I8").Value = 1.053723140879318
$ws.Range("J8").Value = 1.059166699961174
$ws.Range("K8").Value = 1.070806001139087
$ws.Range("L8").Value = 1.064453077893558
$ws.Range("M8").Value = 1.077443960290222
$ws.Range("N8").Value = 1.023431040266357
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.051713609288588
$ws.Range("D9").Value = 1.066248468061423
$ws.Range("E9").Value = 1.059559825235911
$ws.Range("F9").Value = 1.072493279292837
$ws.Range("I9").Value = 1.052921689881474
$ws.Range("J9").Value = 1.057416212238314
$ws.Range("K9").Value = 1.069315429198888
$ws.Range("L9").Value = 1.062647415843033
$ws.Range("M9").Value = 1.075541233142731
$ws.Range("N9").Value = 1.022829915733876
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050010295164916
$ws.Range("D10").Value = 1.064973196837669
$ws.Range("E10").Value = 1.058074192425707
$ws.Range("F10").Value = 1.070946777191243
$ws.Range("I10").Value = 1.052380839531436
$ws.Range("J10").Value = 1.056247092713328
$ws.Range("K10").Value = 1.068318374181673
$ws.Range("L10").Value = 1.061442733107868
$ws.Range("M10").Value = 1.074271992041709
$ws.Range("N10").Value = 1.022427550757014
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049272930936336
$ws.Range("D11").Value = 1.064421006170787
$ws.Range("E11").Value = 1.057431447803341
$ws.Range("F11").Value = 1.070277705707565
$ws.Range("I11").Value = 1.052145093221428
$ws.Range("J11").Value = 1.055740346178991
$ws.Range("K11").Value = 1.067885850181401
$ws.Range("L11").Value = 1.060920878301701
$ws.Range("M11").Value = 1.073722220370754
$ws.Range("N11").Value = 1.022252940489828
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.048999067710976
$ws.Range("D12").Value = 1.064215899998203
$ws.Range("E12").Value = 1.057192785736656
$ws.Range("F12").Value = 1.070029269276601
$ws.Range("I12").Value = 1.052057293033422
$ws.Range("J12").Value = 1.055552041205108
$ws.Range("K12").Value = 1.067725072682296
$ws.Range("L12").Value = 1.060727005283041
$ws.Range("M12").Value = 1.073517983526901
$ws.Range("N12").Value = 1.022188024907119
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049057811105802
$ws.Range("D13").Value = 1.06425989586877
$ws.Range("E13").Value = 1.057243975851424
$ws.Range("F13").Value = 1.070082555815647
$ws.Range("I13").Value = 1.052076137038038
$ws.Range("J13").Value = 1.055592436736093
$ws.Range("K13").Value = 1.067759565374772
$ws.Range("L13").Value = 1.060768593192342
$ws.Range("M13").Value = 1.073561794249958
$ws.Range("N13").Value = 1.022201952123275
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049250292776269
$ws.Range("D14").Value = 1.064404051986352
$ws.Range("E14").Value = 1.057411718250787
$ws.Range("F14").Value = 1.070257168108478
$ws.Range("I14").Value = 1.052137840393831
$ws.Range("J14").Value = 1.055724782397415
$ws.Range("K14").Value = 1.067872562685755
$ws.Range("L14").Value = 1.060904853368642
$ws.Range("M14").Value = 1.073705338642785
$ws.Range("N14").Value = 1.02224757572119
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.049368890601653
$ws.Range("D15").Value = 1.064492871611632
$ws.Range("E15").Value = 1.057515080736731
$ws.Range("F15").Value = 1.070364763975775
$ws.Range("I15").Value = 1.052175826926444
$ws.Range("J15").Value = 1.055806314753984
$ws.Range("K15").Value = 1.067942168320198
$ws.Range("L15").Value = 1.060988803400671
$ws.Range("M15").Value = 1.073793777485765
$ws.Range("N15").Value = 1.022275678301241
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.050059235378371
$ws.Range("D16").Value = 1.065009844170721
$ws.Range("E16").Value = 1.05811686074734
$ws.Range("F16").Value = 1.070991193379851
$ws.Range("I16").Value = 1.052396452488648
$ws.Range("J16").Value = 1.056280713055401
$ws.Range("K16").Value = 1.068347062697434
$ws.Range("L16").Value = 1.061477362276148
$ws.Range("M16").Value = 1.074308474758261
$ws.Range("N16").Value = 1.022439130990088
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.050492318679986
$ws.Range("D17").Value = 1.065334130477771
$ws.Range("E17").Value = 1.058494487246914
$ws.Range("F17").Value = 1.071384289884111
$ws.Range("I17").Value = 1.052534428908569
$ws.Range("J17").Value = 1.056578153934863
$ws.Range("K17").Value = 1.068600830219852
$ws.Range("L17").Value = 1.061783763676821
$ws.Range("M17").Value = 1.074631282173802
$ws.Range("N17").Value = 1.022541557900529
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.050744946450169
$ws.Range("D18").Value = 1.065523282144015
$ws.Range("E18").Value = 1.058714802891624
$ws.Range("F18").Value = 1.071613631802892
$ws.Range("I18").Value = 1.052614758226656
$ws.Range("J18").Value = 1.056751596961416
$ws.Range("K18").Value = 1.068748772031169
$ws.Range("L18").Value = 1.061962461275203
$ws.Range("M18").Value = 1.074819552716405
$ws.Range("N18").Value = 1.022601264769919
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.050831088998209
$ws.Range("D19").Value = 1.065587778120059
$ws.Range("E19").Value = 1.058789933732871
$ws.Range("F19").Value = 1.071691840855966
$ws.Range("I19").Value = 1.052642123005568
$ws.Range("J19").Value = 1.056810728208658
$ws.Range("K19").Value = 1.068799203383792
$ws.Range("L19").Value = 1.062023388949213
$ws.Range("M19").Value = 1.074883745122455
$ws.Range("N19").Value = 1.022621616988888
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05044585112382
$ws.Range("D20").Value = 1.065299337524125
$ws.Range("E20").Value = 1.058453966066222
$ws.Range("F20").Value = 1.071342108619454
$ws.Range("I20").Value = 1.052519640858671
$ws.Range("J20").Value = 1.056546246433813
$ws.Range("K20").Value = 1.068573611276266
$ws.Range("L20").Value = 1.061750891894112
$ws.Range("M20").Value = 1.074596649808332
$ws.Range("N20").Value = 1.022530572290426
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049193610990628
$ws.Range("D21").Value = 1.064361601542773
$ws.Range("E21").Value = 1.057362320032907
$ws.Range("F21").Value = 1.070205746750255
$ws.Range("I21").Value = 1.052119676739206
$ws.Range("J21").Value = 1.055685812018336
$ws.Range("K21").Value = 1.067839291072247
$ws.Range("L21").Value = 1.0608647290511
$ws.Range("M21").Value = 1.07366306914548
$ws.Range("N21").Value = 1.02223414229695
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048406431954198
$ws.Range("D22").Value = 1.06377202168742
$ws.Range("E22").Value = 1.056676432790283
$ws.Range("F22").Value = 1.069491771577307
$ws.Range("I22").Value = 1.051866852416181
$ws.Range("J22").Value = 1.055144378653334
$ws.Range("K22").Value = 1.067376907441597
$ws.Range("L22").Value = 1.060307372919501
$ws.Range("M22").Value = 1.073075932831482
$ws.Range("N22").Value = 1.022047432036162
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.048823716076262
$ws.Range("D23").Value = 1.064084567823439
$ws.Range("E23").Value = 1.057039989621456
$ws.Range("F23").Value = 1.069870215786725
$ws.Range("I23").Value = 1.052001007419447
$ws.Range("J23").Value = 1.055431444858215
$ws.Range("K23").Value = 1.067622090823836
$ws.Range("L23").Value = 1.060602855986246
$ws.Range("M23").Value = 1.073387199737502
$ws.Range("N23").Value = 1.022146442214166
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.050466847757611
$ws.Range("D24").Value = 1.065315058959265
$ws.Range("E24").Value = 1.058472275681258
$ws.Range("F24").Value = 1.071361168346084
$ws.Range("I24").Value = 1.052526323405262
$ws.Range("J24").Value = 1.056560664211997
$ws.Range("K24").Value = 1.068585910581528
$ws.Range("L24").Value = 1.061765745303024
$ws.Range("M24").Value = 1.074612298737725
$ws.Range("N24").Value = 1.022535536329101
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.052374303623594
$ws.Range("D25").Value = 1.066743010807163
$ws.Range("E25").Value = 1.060136422644052
$ws.Range("F25").Value = 1.073093507187809
$ws.Range("I25").Value = 1.053130039234468
$ws.Range("J25").Value = 1.057869130497676
$ws.Range("K25").Value = 1.069701368181774
$ws.Range("L25").Value = 1.063114383583813
$ws.Range("M25").Value = 1.076033268367597
$ws.Range("N25").Value = 1.022985606091618
